# Final update before public release for Feb meeting
#
# grouping_vars.xlsx: "Farm Size" -> "Livestock Area" and
# "Male/Female Household Head" -> "Household Head Gender" (shortName column,
# column C) for the livestock_area and fhh rows respectively, then
# auto-fit the varName/label columns (A/B) to the refreshed content and
# leave the selection on C3.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5 = livestock_area: shortName "Farm Size" -> "Livestock Area"
$ws.Range("C5").Value = "Livestock Area"

# Row 2 = fhh: shortName "Male/Female Household Head" -> "Household Head Gender"
$ws.Range("C2").Value = "Household Head Gender"

# Re-fit the varName (A) and label (B) columns now that content changed
$ws.Columns.Item(1).AutoFit()
$ws.Columns.Item(2).AutoFit()

# Leave the selection where the author left it before saving
[void]$ws.Range("C3").Select()
